$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.319.69'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '''1.679.47'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '''218.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = '''0.5267'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.83%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '''0.2694'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('D9').Value = '''0.06465'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '''21.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('D11').Value = '''0.07512'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '''1.685.11'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').Value = '''0.5802'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').Value = '''0.000008512'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').Value = '''64.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').Value = '''26.335.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '''4.928'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '''190.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '''6.207'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D24').Value = '''145.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '''7.808'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.16%  '
$ws.Range('D26').Value = '''0.1252'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.20%  '
$ws.Range('D27').Value = '''15.79'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = '''0.06497'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('D29').Value = '''1.363'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.55%  '
$ws.Range('D31').Value = '''3.594'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').Value = '''3.589'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('D33').Value = '''1.662'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('D34').Value = '''1.028'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('E35').Value = '  +2.10%  '
$ws.Range('D36').Value = '''2.406'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.88%  '
$ws.Range('E37').Value = '  +3.43%  '
$ws.Range('D38').Value = '''6.339'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.21%  '
$ws.Range('D39').Value = '''1.110.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.75%  '
$ws.Range('D40').Value = '''0.01621'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').Value = '''0.8746'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('D43').Value = '''100.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('D45').Value = '''0.00000000111'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').Value = '''56.87'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = '''8.200'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = '''0.05269'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = '''0.4295'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '''6.074'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.00%  '
